$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes the old row 6 "賣掉/nv/-" and
# everything below it down by one row).
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the "被 / passive / 虛詞" entry.
$ws.Range("A6").Value = "被"
$ws.Range("B6").Value = "passive"
$ws.Range("C6").Value = "虛詞"

# Make row 5 (the row right above the new entry) and the new row 6 a bit
# taller, matching the visual separator look of the source edit.
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5

# Highlight the new row: red font colour plus a red, medium-weight box
# border drawn around the A6:C6 block.
$newRow = $ws.Range("A6:C6")
$newRow.Font.Color = 255
$newRow.BorderAround(1, -4138, 0, 255)

# Leave the same block selected, as in the saved file.
$newRow.Select()
